$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 332.33334
$ws.Range("I12").Value = 332.33334
$ws.Range("K12").Value = 332.33334
$ws.Range("M12").Value = -162.33334

$ws.Range("H58").Value = 3150.5833
$ws.Range("I58").Value = 202.25
$ws.Range("J58").Value = 4624.75
$ws.Range("K58").Value = 606.75
$ws.Range("L58").Value = 13874.25
$ws.Range("M58").Value = -456.75
$ws.Range("N58").Value = -14174.25

$ws.Range("H69").Value = 15499.5
$ws.Range("J69").Value = 18374
$ws.Range("L69").Value = 55122
$ws.Range("N69").Value = -56870

$ws.Range("H72").Value = 15499.5
$ws.Range("J72").Value = 18374
$ws.Range("L72").Value = 165366
$ws.Range("N72").Value = -174102

$ws.Range("H97").Value = 2524
$ws.Range("J97").Value = 2524
$ws.Range("L97").Value = 7572
$ws.Range("N97").Value = -8564

$ws.Range("H112").Value = 2416.2222
$ws.Range("J112").Value = 2470.1177
$ws.Range("L112").Value = 7410.353099999999
$ws.Range("N112").Value = -9626.3531

$ws.Range("H115").Value = 820.75
$ws.Range("I115").Value = 427.66666
$ws.Range("K115").Value = 1282.99998
$ws.Range("M115").Value = 284.0000199999999

$ws.Range("H118").Value = 709.9
$ws.Range("I118").Value = 637.375
$ws.Range("K118").Value = 1912.125
$ws.Range("M118").Value = -255.125

$ws.Range("H132").Value = 6462.6294
$ws.Range("I132").Value = 5779.84
$ws.Range("K132").Value = 17339.52
$ws.Range("M132").Value = -14809.52

$ws.Range("H137").Value = 2739.486
$ws.Range("I137").Value = 1347.3793
$ws.Range("J137").Value = 3678.3489
$ws.Range("K137").Value = 4042.1379
$ws.Range("L137").Value = 11035.0467
$ws.Range("M137").Value = -1492.1379
$ws.Range("N137").Value = -16135.0467

$ws.Range("H138").Value = 3445.7058
$ws.Range("J138").Value = 4565.125
$ws.Range("L138").Value = 13695.375
$ws.Range("N138").Value = -23975.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1087957.6
$ws.Range("I32").Value = 489982.25
$ws.Range("K32").Value = 489982.25
$ws.Range("M32").Value = -489695.25

$ws.Range("H61").Value = 3855.25
$ws.Range("I61").Value = 1696.3334
$ws.Range("J61").Value = 7741.3
$ws.Range("K61").Value = 1696.3334
$ws.Range("L61").Value = 7741.3
$ws.Range("M61").Value = -1484.3334
$ws.Range("N61").Value = -8165.3

$ws.Range("H74").Value = 25005536
$ws.Range("I74").Value = 3517.4443
$ws.Range("J74").Value = 45461736
$ws.Range("K74").Value = 3517.4443
$ws.Range("L74").Value = 45461736
$ws.Range("M74").Value = -2643.4443
$ws.Range("N74").Value = -45463484

$ws.Range("H77").Value = 25005536
$ws.Range("I77").Value = 3517.4443
$ws.Range("J77").Value = 45461736
$ws.Range("K77").Value = 17587.2215
$ws.Range("L77").Value = 227308680
$ws.Range("M77").Value = -13219.2215
$ws.Range("N77").Value = -227317416

$ws.Range("H132").Value = 3547.3794
$ws.Range("I132").Value = 2624.6316
$ws.Range("J132").Value = 5300.6
$ws.Range("K132").Value = 7873.8948
$ws.Range("L132").Value = 15901.8
$ws.Range("M132").Value = -5343.8948
$ws.Range("N132").Value = -20961.8

$ws.Range("H136").Value = 3855.25
$ws.Range("I136").Value = 1696.3334
$ws.Range("J136").Value = 7741.3
$ws.Range("K136").Value = 5089.0002
$ws.Range("L136").Value = 23223.9
$ws.Range("M136").Value = -2539.0002
$ws.Range("N136").Value = -28323.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 133336696
$ws.Range("I94").Value = 148151330
$ws.Range("K94").Value = 148151330
$ws.Range("M94").Value = -148150879

$ws.Range("H105").Value = 10401316
$ws.Range("I105").Value = 477114.94
$ws.Range("K105").Value = 477114.94
$ws.Range("M105").Value = -475367.94

$ws.Range("H134").Value = 2896.087
$ws.Range("I134").Value = 2066
$ws.Range("J134").Value = 3534.6155
$ws.Range("K134").Value = 6198
$ws.Range("L134").Value = 10603.8465
$ws.Range("M134").Value = -3663
$ws.Range("N134").Value = -15673.8465

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2749.25
$ws.Range("I22").Value = 2749.25
$ws.Range("K22").Value = 2749.25
$ws.Range("M22").Value = -2399.25

$ws.Range("H31").Value = 2553504.5
$ws.Range("I31").Value = 1359.6666
$ws.Range("J31").Value = 5003563.5
$ws.Range("K31").Value = 1359.6666
$ws.Range("L31").Value = 5003563.5
$ws.Range("M31").Value = -1064.6666
$ws.Range("N31").Value = -5004153.5

$ws.Range("H34").Value = 2553504.5
$ws.Range("I34").Value = 1359.6666
$ws.Range("J34").Value = 5003563.5
$ws.Range("K34").Value = 1359.6666
$ws.Range("L34").Value = 5003563.5
$ws.Range("M34").Value = -1157.6666
$ws.Range("N34").Value = -5003967.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 232.2381
$ws.Range("I2").Value = 173.42857
$ws.Range("J2").Value = 349.85715
$ws.Range("K2").Value = 1040.57142
$ws.Range("L2").Value = 2099.1429
$ws.Range("M2").Value = -927.57142
$ws.Range("N2").Value = -2325.1429

$ws.Range("H3").Value = 2809.3
$ws.Range("I3").Value = 2565.889
$ws.Range("K3").Value = 7697.667
$ws.Range("M3").Value = -7585.667

$ws.Range("H68").Value = 4551043.5
$ws.Range("J68").Value = 7150278
$ws.Range("L68").Value = 21450834
$ws.Range("N68").Value = -21452456

$ws.Range("H71").Value = 4551043.5
$ws.Range("J71").Value = 7150278
$ws.Range("L71").Value = 64352502
$ws.Range("N71").Value = -64360614

$ws.Range("H111").Value = 723.3333
$ws.Range("I111").Value = 723.3333
$ws.Range("K111").Value = 2169.9999
$ws.Range("M111").Value = 897.0001000000002

$ws.Range("H131").Value = 2188527.5
$ws.Range("I131").Value = 11094
$ws.Range("J131").Value = 14708771
$ws.Range("K131").Value = 33282
$ws.Range("L131").Value = 44126313
$ws.Range("M131").Value = -28242
$ws.Range("N131").Value = -44136393

$ws.Range("H132").Value = 6866.852
$ws.Range("J132").Value = 7596.9287
$ws.Range("L132").Value = 68372.35830000001
$ws.Range("N132").Value = -73432.35830000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 5000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 5000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 5000
$ws.Range("N6").Value = -5226
$ws.Range("M6").ClearContents()

$ws.Range("H16").Value = 5000
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5000
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5000
$ws.Range("N16").Value = -5500
$ws.Range("M16").ClearContents()

$ws.Range("H122").Value = 90915510
$ws.Range("I122").Value = 71434020
$ws.Range("K122").Value = 214302060
$ws.Range("M122").Value = -214299610

$ws.Range("H132").Value = 1884.56
$ws.Range("I132").Value = 1730.7
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 5192.1
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -2662.1
$ws.Range("N132").Value = -12560

$ws.Range("H136").Value = 89176.46000000001
$ws.Range("J136").Value = 89176.46000000001
$ws.Range("L136").Value = 267529.38
$ws.Range("N136").Value = -272629.38

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 133931010
$ws.Range("I22").Value = 14287138
$ws.Range("J22").Value = 333337470
$ws.Range("K22").Value = 14287138
$ws.Range("L22").Value = 333337470
$ws.Range("M22").Value = -14286843
$ws.Range("N22").Value = -333338060

$ws.Range("H27").Value = 133931010
$ws.Range("I27").Value = 14287138
$ws.Range("J27").Value = 333337470
$ws.Range("K27").Value = 14287138
$ws.Range("L27").Value = 333337470
$ws.Range("M27").Value = -14287031
$ws.Range("N27").Value = -333337684

$ws.Range("H46").Value = 5046.4707
$ws.Range("I46").Value = 4599.6665
$ws.Range("J46").Value = 5142.2144
$ws.Range("K46").Value = 4599.6665
$ws.Range("L46").Value = 5142.2144
$ws.Range("M46").Value = -4411.6665
$ws.Range("N46").Value = -5518.2144

$ws.Range("H136").Value = 6883.087
$ws.Range("I136").Value = 5315.75
$ws.Range("J136").Value = 8592.909
$ws.Range("K136").Value = 15947.25
$ws.Range("L136").Value = 25778.727
$ws.Range("M136").Value = -13397.25
$ws.Range("N136").Value = -30878.727

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 2332.6667
$ws.Range("I8").Value = 2332.6667
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 2332.6667
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2192.6667
$ws.Range("N8").ClearContents()

$ws.Range("H28").Value = 4500
$ws.Range("J28").Value = 4500
$ws.Range("L28").Value = 4500
$ws.Range("N28").Value = -5196

$ws.Range("H37").Value = 30000
$ws.Range("I37").Value = 30000
$ws.Range("K37").Value = 30000
$ws.Range("M37").Value = -29797

$ws.Range("H96").Value = 1499.75
$ws.Range("J96").Value = 1500
$ws.Range("L96").Value = 1500
$ws.Range("N96").Value = -4246

$ws.Range("H126").Value = 10749.538
$ws.Range("I126").Value = 14105.111
$ws.Range("J126").Value = 3199.5
$ws.Range("K126").Value = 42315.333
$ws.Range("L126").Value = 9598.5
$ws.Range("M126").Value = -39845.333
$ws.Range("N126").Value = -14538.5

$ws.Range("H136").Value = 13895015
$ws.Range("I136").Value = 15879874
$ws.Range("K136").Value = 47639622
$ws.Range("M136").Value = -47637072
